$d = $word.ActiveDocument

# 1. CL row: F value 53.90 -> 54.01
$d.Content.Find.Execute("53.90", $true, $false, $false, $false, $false,
                         $true, 1, $false, "54.01", 2)

# 2. Evaluation row: F value 9.18 -> 9.19
$d.Content.Find.Execute("9.18", $true, $false, $false, $false, $false,
                         $true, 1, $false, "9.19", 2)

# 3. CL:Evaluation row: F value 1.01 -> 1.03
$d.Content.Find.Execute("1.01", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1.03", 2)

# 4. CL:Evaluation row: p value .321 -> .317
$d.Content.Find.Execute(".321", $true, $false, $false, $false, $false,
                         $true, 1, $false, ".317", 2)

# 5. CL:Evaluation row: 95% CI [0.00, 0.04] -> [0.00, 0.07]
$d.Content.Find.Execute("[0.00, 0.04]", $true, $false, $false, $false, $false,
                         $true, 1, $false, "[0.00, 0.07]", 2)
